# Auto-generated edit script: updates Leve profit-calculation values across
# the ALC, ARM, BSM, CRP, GSM, LTW and WVR sheets (market price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2489.9
$ws.Range("I2").Value = 273.42856
$ws.Range("J2").Value = 7661.6665
$ws.Range("K2").Value = 273.42856
$ws.Range("L2").Value = 7661.6665
$ws.Range("M2").Value = -160.42856
$ws.Range("N2").Value = -7887.6665
$ws.Range("H40").Value = 1945.9512
$ws.Range("I40").Value = 1749.7
$ws.Range("K40").Value = 1749.7
$ws.Range("M40").Value = -1574.7
$ws.Range("H86").Value = 2179.8
$ws.Range("I86").Value = 1766.6666
$ws.Range("J86").Value = 2799.5
$ws.Range("K86").Value = 1766.6666
$ws.Range("L86").Value = 2799.5
$ws.Range("M86").Value = -643.6666
$ws.Range("N86").Value = -5045.5
$ws.Range("H89").Value = 2179.8
$ws.Range("I89").Value = 1766.6666
$ws.Range("J89").Value = 2799.5
$ws.Range("K89").Value = 8833.333000000001
$ws.Range("L89").Value = 13997.5
$ws.Range("M89").Value = -3217.333000000001
$ws.Range("N89").Value = -25229.5
$ws.Range("H94").Value = 1095.25
$ws.Range("I94").Value = 1127
$ws.Range("K94").Value = 1127
$ws.Range("M94").Value = -676
$ws.Range("H107").Value = 415.66666
$ws.Range("I107").Value = 403.33334
$ws.Range("K107").Value = 403.33334
$ws.Range("M107").Value = 1516.66666
$ws.Range("H129").Value = 3190.1428
$ws.Range("I129").Value = 2799
$ws.Range("J129").Value = 3346.6
$ws.Range("K129").Value = 8397
$ws.Range("L129").Value = 10039.8
$ws.Range("M129").Value = -3397
$ws.Range("N129").Value = -20039.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3494.5
$ws.Range("I61").Value = 3193.4
$ws.Range("K61").Value = 3193.4
$ws.Range("M61").Value = -2981.4
$ws.Range("H74").Value = 1771
$ws.Range("I74").Value = 1716.25
$ws.Range("J74").Value = 1990
$ws.Range("K74").Value = 1716.25
$ws.Range("L74").Value = 1990
$ws.Range("M74").Value = -842.25
$ws.Range("N74").Value = -3738
$ws.Range("H77").Value = 1771
$ws.Range("I77").Value = 1716.25
$ws.Range("J77").Value = 1990
$ws.Range("K77").Value = 8581.25
$ws.Range("L77").Value = 9950
$ws.Range("M77").Value = -4213.25
$ws.Range("N77").Value = -18686
$ws.Range("H109").Value = 75000
$ws.Range("J109").Value = 75000
$ws.Range("L109").Value = 75000
$ws.Range("N109").Value = -77774
$ws.Range("H136").Value = 3494.5
$ws.Range("I136").Value = 3193.4
$ws.Range("K136").Value = 9580.200000000001
$ws.Range("M136").Value = -7030.200000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2699.6667
$ws.Range("I86").Value = 2649.5
$ws.Range("K86").Value = 2649.5
$ws.Range("M86").Value = -1526.5
$ws.Range("H89").Value = 2699.6667
$ws.Range("I89").Value = 2649.5
$ws.Range("K89").Value = 13247.5
$ws.Range("M89").Value = -7631.5
$ws.Range("H99").Value = 1500
$ws.Range("J99").Value = 1500
$ws.Range("L99").Value = 1500
$ws.Range("N99").Value = -4496
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 245.75
$ws.Range("I7").Value = 88.625
$ws.Range("K7").Value = 88.625
$ws.Range("M7").Value = 24.375
$ws.Range("H58").Value = 2195.3
$ws.Range("I58").Value = 1918.6
$ws.Range("J58").Value = 2472
$ws.Range("K58").Value = 1918.6
$ws.Range("L58").Value = 2472
$ws.Range("M58").Value = -1715.6
$ws.Range("N58").Value = -2878
$ws.Range("H80").Value = 35000
$ws.Range("J80").Value = 35000
$ws.Range("L80").Value = 35000
$ws.Range("N80").Value = -37246
$ws.Range("H83").Value = 35000
$ws.Range("J83").Value = 35000
$ws.Range("L83").Value = 105000
$ws.Range("N83").Value = -116232
$ws.Range("H92").Value = 29499
$ws.Range("J92").Value = 29499
$ws.Range("L92").Value = 29499
$ws.Range("N92").Value = -34491
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 1439.4
$ws.Range("I122").Value = 1100
$ws.Range("K122").Value = 3300
$ws.Range("M122").Value = -850
$ws.Range("H132").Value = 4662
$ws.Range("I132").Value = 4795
$ws.Range("K132").Value = 14385
$ws.Range("M132").Value = -11855
$ws.Range("H136").Value = 2195.3
$ws.Range("I136").Value = 1918.6
$ws.Range("J136").Value = 2472
$ws.Range("K136").Value = 5755.799999999999
$ws.Range("L136").Value = 7416
$ws.Range("M136").Value = -3205.799999999999
$ws.Range("N136").Value = -12516
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 250750
$ws.Range("I14").Value = 334166.66
$ws.Range("K14").Value = 334166.66
$ws.Range("M14").Value = -333998.66
$ws.Range("H43").Value = 7511.25
$ws.Range("H57").Value = 23055
$ws.Range("I57").Value = 23055
$ws.Range("K57").Value = 23055
$ws.Range("M57").Value = -22235
$ws.Range("H70").Value = 4908
$ws.Range("I70").Value = 4927.5
$ws.Range("J70").Value = 4888.5
$ws.Range("K70").Value = 4927.5
$ws.Range("L70").Value = 4888.5
$ws.Range("M70").Value = -4657.5
$ws.Range("N70").Value = -5428.5
$ws.Range("H73").Value = 4908
$ws.Range("I73").Value = 4927.5
$ws.Range("J73").Value = 4888.5
$ws.Range("K73").Value = 4927.5
$ws.Range("L73").Value = 4888.5
$ws.Range("M73").Value = -3991.5
$ws.Range("N73").Value = -6760.5
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("H113").Value = 999.2
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1507.3143
$ws.Range("J46").Value = 1989.9333
$ws.Range("L46").Value = 1989.9333
$ws.Range("N46").Value = -2365.9333
$ws.Range("H61").Value = 1899.6666
$ws.Range("I61").Value = 1899.6666
$ws.Range("K61").Value = 1899.6666
$ws.Range("M61").Value = -1697.6666
$ws.Range("H113").Value = 1899.6666
$ws.Range("I113").Value = 1899.6666
$ws.Range("K113").Value = 1899.6666
$ws.Range("M113").Value = 270.3334
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").ClearContents()
$ws.Range("N131").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13166.333
$ws.Range("I62").Value = 17916.666
$ws.Range("J62").Value = 3665.6667
$ws.Range("K62").Value = 17916.666
$ws.Range("L62").Value = 3665.6667
$ws.Range("M62").Value = -17292.666
$ws.Range("N62").Value = -4913.6667
$ws.Range("H65").Value = 13166.333
$ws.Range("I65").Value = 17916.666
$ws.Range("J65").Value = 3665.6667
$ws.Range("K65").Value = 89583.33
$ws.Range("L65").Value = 18328.3335
$ws.Range("M65").Value = -86463.33
$ws.Range("N65").Value = -24568.3335
$ws.Range("H81").Value = 1002069.5
$ws.Range("J81").Value = 1669082.9
$ws.Range("L81").Value = 3338165.8
$ws.Range("N81").Value = -3340287.8
$ws.Range("H84").Value = 1002069.5
$ws.Range("J84").Value = 1669082.9
$ws.Range("L84").Value = 16690829
$ws.Range("N84").Value = -16701437
$ws.Range("H100").Value = 9091163
$ws.Range("I100").Value = 11111372
$ws.Range("K100").Value = 22222744
$ws.Range("M100").Value = -22222203
$ws.Range("H107").Value = 3599.8333
$ws.Range("I107").Value = 3700
$ws.Range("J107").Value = 3399.5
$ws.Range("K107").Value = 11100
$ws.Range("L107").Value = 10198.5
$ws.Range("M107").Value = -9180
$ws.Range("N107").Value = -14038.5
$ws.Range("H113").Value = 586.25
$ws.Range("I113").Value = 501
$ws.Range("K113").Value = 1503
$ws.Range("M113").Value = 667
$ws.Range("H126").Value = 2831.7778
$ws.Range("I126").Value = 2017.4
$ws.Range("K126").Value = 6052.200000000001
$ws.Range("M126").Value = -3582.200000000001
$ws.Range("H132").Value = 4489.8
$ws.Range("I132").Value = 4489.8
$ws.Range("K132").Value = 13469.4
$ws.Range("M132").Value = -10939.4
